# adding averages and more checks
# - Training Dashboard: refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I)
#   columns for every data row, and make the header row text white-on-navy.
# - Exam Dashboard: replace the per-row COMMENTS with "date is valid",
#   narrow the COMMENTS column, and make the header row text white-on-navy.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Training Dashboard"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Training Dashboard")

# New "PERIOD TO EXPIRE" values (H3:H26) - last update moved from 08-Sep-2025
# to 16-Sep-2025, so every period-to-expire count drops by 8 days.
$periodToExpire = @{
    3  = 241
    4  = 231
    5  = 239
    6  = 363
    7  = 244
    8  = 219
    9  = 247
    10 = 238
    11 = 491
    12 = 223
    13 = 349
    14 = 254
    15 = 360
    16 = -23
    17 = -103
    18 = -126
    19 = -34
    20 = -34
    21 = 155
    22 = 268
    23 = 313
    24 = 313
    25 = 313
    26 = 348
}

foreach ($row in $periodToExpire.Keys) {
    $ws1.Cells.Item($row, 8).Value = $periodToExpire[$row]
}

# Force the LAST UPDATE column (I) to stay plain text (it is not a real
# date cell) while we overwrite it, otherwise Excel will silently convert
# the "16-Sep-2025" literal into a date serial.
$ws1.Range("I3:I15").NumberFormat = "@"
$ws1.Range("I16:I20").NumberFormat = "@"
$ws1.Range("I21:I26").NumberFormat = "@"

for ($row = 3; $row -le 26; $row++) {
    $ws1.Cells.Item($row, 9).Value = "16-Sep-2025"
}

# Header row (row 2) text turns white so it reads clearly on the navy fill.
$ws1.Range("A2:K2").Font.Color = 16777215

# ---------------------------------------------------------------------------
# Sheet 2: "Exam Dashboard"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

for ($row = 3; $row -le 12; $row++) {
    $ws2.Cells.Item($row, 5).Value = "date is valid"
}

# COMMENTS column no longer needs to be as wide now that every remark is
# the short "date is valid" string.
$ws2.Columns.Item(5).ColumnWidth = 14.1666666666667

# Header row (row 2) text turns white, same treatment as the other sheet.
$ws2.Range("A2:G2").Font.Color = 16777215
